{"js": "// Update the \"feature selection\" paragraph in the final report:\n//   \"if one feature which the F value greater than XXX added...\"\n// becomes\n//   \"if duplicating \"one\" feature which the F value greater than 8.0\n//    (\"shared_receipt_with_poi\", \"salary\", \"total_payments\",\n//    \"total_stock_value\", \"exercised_stock_options\",  \"restricted_stock\")\n//    added...\"\n// The newly added text is colored blue, and the document's \"_GoBack\"\n// bookmark is moved from the (now blank) trailing paragraph onto this\n// new span.\n\nconst body = context.document.body;\n\n// 1) Remove the stale \"_GoBack\" bookmark that currently sits by itself in\n//    an empty paragraph near the end of the document. Use\n//    document.deleteBookmark (not range.delete(), which would also\n//    swallow the now-empty paragraph mark) so the blank paragraph itself\n//    is left in place.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Find the sentence fragment that needs to be expanded.\nconst oldFragment =\n  \"if one feature which the F value greater than XXX added into the \" +\n  \"algorithms, the average accuracy, precision and recall would greatly \" +\n  \"improve.\";\n\nconst found = body.search(oldFragment, { matchCase: true });\nfound.load(\"items\");\nawait context.sync();\n\nif (found.items.length === 0) {\n  throw new Error(\"Could not locate the target sentence to update.\");\n}\n\nconst target = found.items[0];\n\nconst newFragment =\n  \"if duplicating \\u201Cone\\u201D feature which the F value greater than \" +\n  \"8.0 (\\u201Cshared_receipt_with_poi\\u201D, \\u201Csalary\\u201D, \" +\n  \"\\u201Ctotal_payments\\u201D, \\u201Ctotal_stock_value\\u201D, \" +\n  \"\\u201Cexercised_stock_options\\u201D,  \\u201Crestricted_stock\\u201D) \" +\n  \"added into the algorithms, the average accuracy, precision and recall \" +\n  \"would greatly improve.\";\n\n// 3) Swap in the new wording (this keeps the surrounding run formatting\n//    intact on either side of the replaced range).\ntarget.insertText(newFragment, Word.InsertLocation.replace);\nawait context.sync();\n\n// 4) Re-find the whole updated sentence -- including the leading\n//    \"In addition, \" that now also gets the blue highlight -- color it\n//    blue, and wrap it with a new \"_GoBack\" bookmark (matching where Word\n//    last left the cursor after the edit).\nconst blueSentence = \"In addition, \" + newFragment;\nconst inserted = body.search(blueSentence, { matchCase: true });\ninserted.load(\"items\");\nawait context.sync();\n\nconst newRange = inserted.items[0];\nnewRange.font.color = \"#0000FF\";\nnewRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Update the \"feature selection\" paragraph in the final report:\n#   \"if one feature which the F value greater than XXX added...\"\n# becomes\n#   \"if duplicating \"one\" feature which the F value greater than 8.0\n#    (\"shared_receipt_with_poi\", \"salary\", \"total_payments\",\n#    \"total_stock_value\", \"exercised_stock_options\",  \"restricted_stock\")\n#    added...\"\n# The newly added text is colored blue, and the document's \"_GoBack\"\n# bookmark is moved from the (now blank) trailing paragraph onto this\n# new span.\n\n$d = $word.ActiveDocument\n\n# 1) Remove the stale \"_GoBack\" bookmark that currently sits by itself in\n#    an empty paragraph near the end of the document.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2) Replace the old sentence fragment with the expanded wording. Curly\n#    quotes are built from character codes to avoid encoding issues.\n$quoteOpen  = [char]8220\n$quoteClose = [char]8221\n\n$oldText = \"if one feature which the F value greater than XXX added into the algorithms, the average accuracy, precision and recall would greatly improve.\"\n\n$newText = \"if duplicating \" + $quoteOpen + \"one\" + $quoteClose + `\n    \" feature which the F value greater than 8.0 (\" + $quoteOpen + \"shared_receipt_with_poi\" + $quoteClose + `\n    \", \" + $quoteOpen + \"salary\" + $quoteClose + `\n    \", \" + $quoteOpen + \"total_payments\" + $quoteClose + `\n    \", \" + $quoteOpen + \"total_stock_value\" + $quoteClose + `\n    \", \" + $quoteOpen + \"exercised_stock_options\" + $quoteClose + `\n    \",  \" + $quoteOpen + \"restricted_stock\" + $quoteClose + `\n    \") added into the algorithms, the average accuracy, precision and recall would greatly improve.\"\n\n$findRange = $d.Content\n$findRange.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n\n# 3) Color the whole updated sentence -- including the leading\n#    \"In addition, \" that now also gets the blue highlight -- blue\n#    (wdColorBlue = 16711680 -> OOXML w:color \"0000FF\").\n$blueText = \"In addition, \" + $newText\n$blueRange = $d.Content\n$blueRange.Find.Execute($blueText) | Out-Null\n$blueRange.Font.Color = 16711680\n\n# 4) Wrap the freshly-colored span with a new \"_GoBack\" bookmark (matching\n#    where Word last left the cursor after the edit).\n$d.Bookmarks.Add(\"_GoBack\", $blueRange)\n"}
